# Conserto do erro com o rotulo da coluna 2050 nas tabelas e retirada das linhas com total das tabelas

$wb = $excel.ActiveWorkbook

# Sheets with "year" style headers (2015 / 2030 / 2040 / 2050) and a Total row to remove
$yearSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)"
)

foreach ($name in $yearSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("E1").Value = "'2050"
    $ws.Rows.Item(13).Delete()
}

# Sheet with "period" style headers (2015 / 2015-2030 / 2031-2040 / 2041-2050) and a Total row to remove
$ws = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws.Range("E1").Value = "2041-2050"
$ws.Rows.Item(13).Delete()

# Sheet with "year" style headers but no Total row
$ws = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
$ws.Range("E1").Value = "'2050"

# Sheet with a Total row to remove (no year headers)
$ws = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws.Rows.Item(4).Delete()
